# The deck has two identical "Postgresql - Settings and Other Sundries"
# slides back-to-back (slide 6 and slide 7). Remove the duplicate at
# position 7; PowerPoint renumbers/shifts every following slide up by
# one when the slide is deleted, which is exactly what the target
# presentation reflects (the "Logging and Getting Information" slide
# moves into slot 7, "Backup and Restore" into slot 8, "Client Tools"
# into slot 9, "Appendix" into slot 10, and so on).

$p = $ppt.ActivePresentation
$p.Slides.Item(7).Delete()
